$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row (row 38) with the new question entry
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "Number of Connected Components in an Undirected Graph"
$ws.Range("D38").Value = "Graph/UF"
$ws.Range("E38").Value = "medium"
$ws.Range("F38").Value = "leetcode 323"

# Match styling of row 37 (centered numbers for A, left for B, general for D/E/F)
$ws.Range("A38").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B38").HorizontalAlignment = -4131  # xlLeft
$ws.Range("D38").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E38").HorizontalAlignment = -4108  # xlCenter
$ws.Range("F38").HorizontalAlignment = -4108  # xlCenter

# Scroll/selection state to match the diff
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("F38").Select()
